$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells for new columns I and J (copy formatting from H1's header style)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data rows 2..13: I = 1 (constant), J = same value as H
for ($r = 2; $r -le 13; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value()
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
